$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 18 de Agosto de 2020 a las 19:29"

# Row 4: Estados Unidos
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 5624297
$ws.Cells.Item(4, 3).Value = 12270
$ws.Cells.Item(4, 4).Value = 2975698
$ws.Cells.Item(4, 5).Value = 2474439
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 444
$ws.Cells.Item(4, 8).Value = 174160

# Row 6: India
$ws.Cells.Item(6, 1).Value = "India"
$ws.Cells.Item(6, 2).Value = 2762591
$ws.Cells.Item(6, 3).Value = 60987
$ws.Cells.Item(6, 4).Value = 2034680
$ws.Cells.Item(6, 5).Value = 674906
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 1080
$ws.Cells.Item(6, 8).Value = 53005

# Row 12: Chile
$ws.Cells.Item(12, 1).Value = "Chile"
$ws.Cells.Item(12, 2).Value = 388855
$ws.Cells.Item(12, 3).Value = 1353
$ws.Cells.Item(12, 4).Value = 362440
$ws.Cells.Item(12, 5).Value = 15869
$ws.Cells.Item(12, 6).Value = 0
$ws.Cells.Item(12, 7).Value = 33
$ws.Cells.Item(12, 8).Value = 10546

# Row 21: Turquia
$ws.Cells.Item(21, 1).Value = "Turquia"
$ws.Cells.Item(21, 2).Value = 251805
$ws.Cells.Item(21, 3).Value = 1263
$ws.Cells.Item(21, 4).Value = 232913
$ws.Cells.Item(21, 5).Value = 12876
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(21, 7).Value = 20
$ws.Cells.Item(21, 8).Value = 6016

# Row 22: Alemania
$ws.Cells.Item(22, 1).Value = "Alemania"
$ws.Cells.Item(22, 2).Value = 227791
$ws.Cells.Item(22, 3).Value = 1105
$ws.Cells.Item(22, 4).Value = 202900
$ws.Cells.Item(22, 5).Value = 15587
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(22, 7).Value = 8
$ws.Cells.Item(22, 8).Value = 9304

# Row 23: Francia
$ws.Cells.Item(23, 1).Value = "Francia"
$ws.Cells.Item(23, 2).Value = 221267
$ws.Cells.Item(23, 3).Value = 2238
$ws.Cells.Item(23, 4).Value = 84065
$ws.Cells.Item(23, 5).Value = 106773
$ws.Cells.Item(23, 6).Value = 0
$ws.Cells.Item(23, 7).Value = 0
$ws.Cells.Item(23, 8).Value = 30429

# Row 24: Irak
$ws.Cells.Item(24, 1).Value = "Irak"
$ws.Cells.Item(24, 2).Value = 184709
$ws.Cells.Item(24, 3).Value = 4576
$ws.Cells.Item(24, 4).Value = 131840
$ws.Cells.Item(24, 5).Value = 46833
$ws.Cells.Item(24, 6).Value = 0
$ws.Cells.Item(24, 7).Value = 82
$ws.Cells.Item(24, 8).Value = 6036

# Row 33: Israel
$ws.Cells.Item(33, 1).Value = "Israel"
$ws.Cells.Item(33, 2).Value = 96093
$ws.Cells.Item(33, 3).Value = 1342
$ws.Cells.Item(33, 4).Value = 71971
$ws.Cells.Item(33, 5).Value = 23417
$ws.Cells.Item(33, 6).Value = 0
$ws.Cells.Item(33, 7).Value = 13
$ws.Cells.Item(33, 8).Value = 705

# Row 44: Emiratos Arabes Unidos
$ws.Cells.Item(44, 1).Value = "Emiratos Arabes Unidos"
$ws.Cells.Item(44, 2).Value = 64906
$ws.Cells.Item(44, 3).Value = 365
$ws.Cells.Item(44, 4).Value = 57909
$ws.Cells.Item(44, 5).Value = 6631
$ws.Cells.Item(44, 6).Value = 0
$ws.Cells.Item(44, 7).Value = 2
$ws.Cells.Item(44, 8).Value = 366

# Row 54: Marruecos
$ws.Cells.Item(54, 1).Value = "Marruecos"
$ws.Cells.Item(54, 2).Value = 44803
$ws.Cells.Item(54, 3).Value = 1245
$ws.Cells.Item(54, 4).Value = 31002
$ws.Cells.Item(54, 5).Value = 13087
$ws.Cells.Item(54, 6).Value = 0
$ws.Cells.Item(54, 7).Value = 33
$ws.Cells.Item(54, 8).Value = 714

# Row 55: Ghana
$ws.Cells.Item(55, 1).Value = "Ghana"
$ws.Cells.Item(55, 2).Value = 42993
$ws.Cells.Item(55, 3).Value = 340
$ws.Cells.Item(55, 4).Value = 40796
$ws.Cells.Item(55, 5).Value = 1949
$ws.Cells.Item(55, 6).Value = 0
$ws.Cells.Item(55, 7).Value = 9
$ws.Cells.Item(55, 8).Value = 248

# Row 74: Chequia
$ws.Cells.Item(74, 1).Value = "Chequia"
$ws.Cells.Item(74, 2).Value = 20393
$ws.Cells.Item(74, 3).Value = 191
$ws.Cells.Item(74, 4).Value = 15146
$ws.Cells.Item(74, 5).Value = 4846
$ws.Cells.Item(74, 6).Value = 0
$ws.Cells.Item(74, 7).Value = 2
$ws.Cells.Item(74, 8).Value = 401

# Row 89: Libano
$ws.Cells.Item(89, 1).Value = "Libano"
$ws.Cells.Item(89, 2).Value = 9758
$ws.Cells.Item(89, 3).Value = 421
$ws.Cells.Item(89, 4).Value = 2852
$ws.Cells.Item(89, 5).Value = 6799
$ws.Cells.Item(89, 6).Value = 0
$ws.Cells.Item(89, 7).Value = 2
$ws.Cells.Item(89, 8).Value = 107

# Row 90: Consejo Danes para los Refugiados
$ws.Cells.Item(90, 1).Value = "Consejo Danes para los Refugiados"
$ws.Cells.Item(90, 2).Value = 9721
$ws.Cells.Item(90, 3).Value = 15
$ws.Cells.Item(90, 4).Value = 8882
$ws.Cells.Item(90, 5).Value = 596
$ws.Cells.Item(90, 6).Value = 0
$ws.Cells.Item(90, 7).Value = 0
$ws.Cells.Item(90, 8).Value = 243

# Row 100: Grecia
$ws.Cells.Item(100, 1).Value = "Grecia"
$ws.Cells.Item(100, 2).Value = 7472
$ws.Cells.Item(100, 3).Value = 250
$ws.Cells.Item(100, 4).Value = 3804
$ws.Cells.Item(100, 5).Value = 3436
$ws.Cells.Item(100, 6).Value = 0
$ws.Cells.Item(100, 7).Value = 2
$ws.Cells.Item(100, 8).Value = 232

# Row 101: Luxemburgo
$ws.Cells.Item(101, 1).Value = "Luxemburgo"
$ws.Cells.Item(101, 2).Value = 7469
$ws.Cells.Item(101, 3).Value = 0
$ws.Cells.Item(101, 4).Value = 6739
$ws.Cells.Item(101, 5).Value = 606
$ws.Cells.Item(101, 6).Value = 0
$ws.Cells.Item(101, 7).Value = 0
$ws.Cells.Item(101, 8).Value = 124

# Row 131: Tunez
$ws.Cells.Item(131, 1).Value = "Tunez"
$ws.Cells.Item(131, 2).Value = 2314
$ws.Cells.Item(131, 3).Value = 129
$ws.Cells.Item(131, 4).Value = 1370
$ws.Cells.Item(131, 5).Value = 887
$ws.Cells.Item(131, 6).Value = 0
$ws.Cells.Item(131, 7).Value = 1
$ws.Cells.Item(131, 8).Value = 57

# Row 132: Estonia
$ws.Cells.Item(132, 1).Value = "Estonia"
$ws.Cells.Item(132, 2).Value = 2200
$ws.Cells.Item(132, 3).Value = 8
$ws.Cells.Item(132, 4).Value = 1990
$ws.Cells.Item(132, 5).Value = 147
$ws.Cells.Item(132, 6).Value = 0
$ws.Cells.Item(132, 7).Value = 0
$ws.Cells.Item(132, 8).Value = 63

# Row 142: Uganda
$ws.Cells.Item(142, 1).Value = "Uganda"
$ws.Cells.Item(142, 2).Value = 1603
$ws.Cells.Item(142, 3).Value = 43
$ws.Cells.Item(142, 4).Value = 1165
$ws.Cells.Item(142, 5).Value = 423
$ws.Cells.Item(142, 6).Value = 0
$ws.Cells.Item(142, 7).Value = 0
$ws.Cells.Item(142, 8).Value = 15

# Row 144: Jordania
$ws.Cells.Item(144, 1).Value = "Jordania"
$ws.Cells.Item(144, 2).Value = 1438
$ws.Cells.Item(144, 3).Value = 40
$ws.Cells.Item(144, 4).Value = 1243
$ws.Cells.Item(144, 5).Value = 184
$ws.Cells.Item(144, 6).Value = 0
$ws.Cells.Item(144, 7).Value = 0
$ws.Cells.Item(144, 8).Value = 11

# Row 213: Montserrat
$ws.Cells.Item(213, 1).Value = "Montserrat"
$ws.Cells.Item(213, 2).Value = 13
$ws.Cells.Item(213, 3).Value = 0
$ws.Cells.Item(213, 4).Value = 12
$ws.Cells.Item(213, 5).Value = 0
$ws.Cells.Item(213, 6).Value = 0
$ws.Cells.Item(213, 7).Value = 0
$ws.Cells.Item(213, 8).Value = 1

# Row 214: Islas Malvinas
$ws.Cells.Item(214, 1).Value = "Islas Malvinas"
$ws.Cells.Item(214, 2).Value = 13
$ws.Cells.Item(214, 3).Value = 0
$ws.Cells.Item(214, 4).Value = 13
$ws.Cells.Item(214, 5).Value = 0
$ws.Cells.Item(214, 6).Value = 0
$ws.Cells.Item(214, 7).Value = 0
$ws.Cells.Item(214, 8).Value = 0

